$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store text exactly as given, without leaving behind
# a residual custom number-format style (matches original "shared string / no style" look).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.469.11"
Set-TextValue $ws.Range("E2") "  +0.49%  "

Set-TextValue $ws.Range("D3") "2.103.00"
Set-TextValue $ws.Range("E3") "  +0.88%  "

Set-TextValue $ws.Range("D4") "1.009"
Set-TextValue $ws.Range("E4") "  +1.01%  "

Set-TextValue $ws.Range("D5") "334.14"
Set-TextValue $ws.Range("E5") "  +1.68%  "

Set-TextValue $ws.Range("D6") "1.006"
Set-TextValue $ws.Range("E6") "  +0.77%  "

Set-TextValue $ws.Range("D7") "0.5215"
Set-TextValue $ws.Range("E7") "  -0.18%  "

Set-TextValue $ws.Range("D8") "0.4517"
Set-TextValue $ws.Range("E8") "  +4.55%  "

Set-TextValue $ws.Range("D9") "53.78"
Set-TextValue $ws.Range("E9") "  +15.05%  "

Set-TextValue $ws.Range("D10") "0.08909"
Set-TextValue $ws.Range("E10") "  +0.80%  "

Set-TextValue $ws.Range("D11") "1.179"
Set-TextValue $ws.Range("E11") "  +1.39%  "

Set-TextValue $ws.Range("D12") "24.06"
Set-TextValue $ws.Range("E12") "  -1.85%  "

Set-TextValue $ws.Range("D13") "2.107.52"
Set-TextValue $ws.Range("E13") "  +1.08%  "

Set-TextValue $ws.Range("D14") "6.809"
Set-TextValue $ws.Range("E14") "  +0.92%  "

Set-TextValue $ws.Range("D15") "8.003"
Set-TextValue $ws.Range("E15") "  +3.96%  "

Set-TextValue $ws.Range("D16") "96.61"
Set-TextValue $ws.Range("E16") "  +1.06%  "

Set-TextValue $ws.Range("B17") "BinanceUSD"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D17") "1.008"
Set-TextValue $ws.Range("E17") "  +0.82%  "

Set-TextValue $ws.Range("B18") "ShibaInu"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.00001140"
Set-TextValue $ws.Range("E18") "  +1.23%  "

Set-TextValue $ws.Range("D19") "0.06651"
Set-TextValue $ws.Range("E19") "  +0.29%  "

Set-TextValue $ws.Range("D20") "19.19"
Set-TextValue $ws.Range("E20") "  +1.41%  "

Set-TextValue $ws.Range("D21") "1.006"
Set-TextValue $ws.Range("E21") "  +0.70%  "

Set-TextValue $ws.Range("D22") "6.317"
Set-TextValue $ws.Range("E22") "  -0.10%  "

Set-TextValue $ws.Range("D23") "30.529.92"
Set-TextValue $ws.Range("E23") "  +0.52%  "

Set-TextValue $ws.Range("D24") "12.44"
Set-TextValue $ws.Range("E24") "  +0.48%  "

Set-TextValue $ws.Range("E25") "  +2.44%  "

Set-TextValue $ws.Range("D26") "2.353.07"
Set-TextValue $ws.Range("E26") "  +1.04%  "

Set-TextValue $ws.Range("D27") "22.18"
Set-TextValue $ws.Range("E27") "  -0.91%  "

Set-TextValue $ws.Range("D28") "162.69"
Set-TextValue $ws.Range("E28") "  +0.50%  "

Set-TextValue $ws.Range("D29") "2.518"
Set-TextValue $ws.Range("E29") "  -2.90%  "

Set-TextValue $ws.Range("D30") "133.25"
Set-TextValue $ws.Range("E30") "  +1.21%  "

Set-TextValue $ws.Range("D31") "1.205"
Set-TextValue $ws.Range("E31") "  +0.57%  "

Set-TextValue $ws.Range("D32") "0.1070"
Set-TextValue $ws.Range("E32") "  -0.04%  "

Set-TextValue $ws.Range("B33") "Filecoin"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "6.398"
Set-TextValue $ws.Range("E33") "  +3.41%  "

Set-TextValue $ws.Range("B34") "ARBITRUM"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D34") "1.628"
Set-TextValue $ws.Range("E34") "  -2.28%  "

Set-TextValue $ws.Range("D35") "3.949"

Set-TextValue $ws.Range("D36") "10.40"
Set-TextValue $ws.Range("E36") "  +3.87%  "

Set-TextValue $ws.Range("D37") "5.810"
Set-TextValue $ws.Range("E37") "  +6.38%  "

Set-TextValue $ws.Range("D38") "0.02579"
Set-TextValue $ws.Range("E38") "  +0.26%  "

Set-TextValue $ws.Range("D39") "0.06835"
Set-TextValue $ws.Range("E39") "  +2.35%  "

Set-TextValue $ws.Range("D40") "0.2290"
Set-TextValue $ws.Range("E40") "  +1.07%  "

Set-TextValue $ws.Range("D41") "12.66"
Set-TextValue $ws.Range("E41") "  -0.31%  "

Set-TextValue $ws.Range("D42") "0.6851"
Set-TextValue $ws.Range("E42") "  +0.23%  "

Set-TextValue $ws.Range("D43") "1.248"
Set-TextValue $ws.Range("E43") "  +0.26%  "

Set-TextValue $ws.Range("B44") "EnergySwap"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "14.07"
Set-TextValue $ws.Range("E44") "  -0.02%  "

Set-TextValue $ws.Range("B45") "NEARProtocol"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D45") "2.308"
Set-TextValue $ws.Range("E45") "  +4.45%  "

Set-TextValue $ws.Range("D46") "0.6350"
Set-TextValue $ws.Range("E46") "  -0.63%  "

Set-TextValue $ws.Range("D47") "3.666"
Set-TextValue $ws.Range("E47") "  +1.61%  "

Set-TextValue $ws.Range("D48") "0.00000000351"
Set-TextValue $ws.Range("E48") "  +22.98%  "

Set-TextValue $ws.Range("D49") "1.247"
Set-TextValue $ws.Range("E49") "  -0.28%  "

Set-TextValue $ws.Range("B50") "WEMIXTOKEN"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D50") "1.207"
Set-TextValue $ws.Range("E50") "  +1.29%  "

Set-TextValue $ws.Range("B51") "Aave"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "83.17"
Set-TextValue $ws.Range("E51") "  +1.79%  "
